$d = $word.ActiveDocument

# The last paragraph in the document currently reads "Fragen" and also
# carries the (hidden) "_GoBack" bookmark at its end. We need to:
#   1) keep the "Fragen" paragraph as-is but drop the bookmark from it
#   2) add a brand new sub-bullet paragraph after it ("Github Problem
#      mit mergen") that now owns the "_GoBack" bookmark near "mergen"
#      plus the matching spell-check proofErr markers.
#
# We do this by replacing the whole last paragraph (identified via the
# Paragraphs collection, independent of any absolute character offsets)
# with a WordprocessingML package fragment describing both paragraphs
# in their final form.

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p w:rsidR="005B597B" w:rsidRPr="005B597B" w:rsidRDefault="00ED693C" w:rsidP="005B597B">' +
'<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr>' +
'<w:r w:rsidRPr="00BE3F73"><w:rPr><w:b/></w:rPr><w:t>Fragen</w:t></w:r>' +
'</w:p>' +
'<w:p>' +
'<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
'<w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
'<w:r><w:t xml:space="preserve"> Problem mit </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/><w:r><w:t>mergen</w:t></w:r>' +
'<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
'<w:proofErr w:type="spellEnd"/>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$target.InsertXML($xml)
